$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5173
$ws1.Range("F7").Value = 61
$ws1.Range("F10").Value = 4
$ws1.Range("F11").Value = 62

# Sheet "全部类型" (All types) - same events repeated, update matching rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 5173
$ws4.Range("F11").Value = 61
$ws4.Range("F15").Value = 4
$ws4.Range("F16").Value = 62
